$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1787878787878788
$ws.Range("C2").Value = 0.5636363636363636
$ws.Range("J2").Value = 0.01212121212121212
$ws.Range("P2").Value = 0.1303030303030303
$ws.Range("S2").Value = 0.1151515151515152
$ws.Range("C3").Value = 0.03157894736842105
$ws.Range("J3").Value = 0.01578947368421053
$ws.Range("P3").Value = 0.7368421052631579
$ws.Range("S3").Value = 0.2157894736842105
$ws.Range("B6").Value = 0.07086614173228346
$ws.Range("D6").Value = 0.007874015748031496
$ws.Range("F6").Value = 0.09842519685039371
$ws.Range("J6").Value = 0.2086614173228346
$ws.Range("O6").Value = 0.03149606299212598
$ws.Range("Q6").Value = 0.1850393700787402
$ws.Range("R6").Value = 0.07874015748031496
$ws.Range("S6").Value = 0.3188976377952756
$ws.Range("B7").Value = 0.1157407407407407
$ws.Range("D7").Value = 0.02777777777777778
$ws.Range("F7").Value = 0.07870370370370371
$ws.Range("J7").Value = 0.1064814814814815
$ws.Range("O7").Value = 0.02777777777777778
$ws.Range("Q7").Value = 0.2222222222222222
$ws.Range("R7").Value = 0.05092592592592592
$ws.Range("S7").Value = 0.3703703703703703
$ws.Range("B8").Value = 0.118942731277533
$ws.Range("D8").Value = 0.01762114537444934
$ws.Range("F8").Value = 0.05506607929515418
$ws.Range("J8").Value = 0.09030837004405286
$ws.Range("O8").Value = 0.02863436123348018
$ws.Range("Q8").Value = 0.2180616740088106
$ws.Range("R8").Value = 0.06167400881057269
$ws.Range("S8").Value = 0.4096916299559472
$ws.Range("B9").Value = 0.12
$ws.Range("D9").Value = 0.03555555555555556
$ws.Range("F9").Value = 0.1155555555555556
$ws.Range("J9").Value = 0.09333333333333334
$ws.Range("O9").Value = 0.004444444444444444
$ws.Range("Q9").Value = 0.1911111111111111
$ws.Range("R9").Value = 0.08444444444444445
$ws.Range("S9").Value = 0.3555555555555556
$ws.Range("B10").Value = 0.111376292760541
$ws.Range("D10").Value = 0.02784407319013524
$ws.Range("E10").Value = 0.0007955449482895784
$ws.Range("F10").Value = 0.07716785998408911
$ws.Range("J10").Value = 0.1058074781225139
$ws.Range("O10").Value = 0.02307080350039777
$ws.Range("Q10").Value = 0.2219570405727924
$ws.Range("R10").Value = 0.05966587112171837
$ws.Range("S10").Value = 0.3723150357995227
$ws.Range("G11").Value = 0.0970873786407767
$ws.Range("J11").Value = 0.06796116504854369
$ws.Range("K11").Value = 0.1618122977346278
$ws.Range("L11").Value = 0.6601941747572816
$ws.Range("S11").Value = 0.01294498381877023
$ws.Range("G12").Value = 0.7649769585253456
$ws.Range("J12").Value = 0.1566820276497696
$ws.Range("L12").Value = 0.04608294930875576
$ws.Range("S12").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.01538461538461539
$ws.Range("H15").Value = 0.1884615384615385
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3076923076923077
$ws.Range("K15").Value = 0.05
$ws.Range("M15").Value = 0.003846153846153846
$ws.Range("O15").Value = 0.06923076923076923
$ws.Range("S15").Value = 0.2884615384615384
$ws.Range("F16").Value = 0.01923076923076923
$ws.Range("H16").Value = 0.1730769230769231
$ws.Range("I16").Value = 0.0673076923076923
$ws.Range("J16").Value = 0.3894230769230769
$ws.Range("K16").Value = 0.1346153846153846
$ws.Range("M16").Value = 0.004807692307692308
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.08653846153846154
$ws.Range("S16").Value = 0.1201923076923077
$ws.Range("F17").Value = 0.01174168297455969
$ws.Range("H17").Value = 0.1741682974559687
$ws.Range("I17").Value = 0.1095890410958904
$ws.Range("J17").Value = 0.4246575342465753
$ws.Range("K17").Value = 0.111545988258317
$ws.Range("M17").Value = 0.01761252446183953
$ws.Range("N17").Value = 0.007827788649706457
$ws.Range("O17").Value = 0.06262230919765166
$ws.Range("S17").Value = 0.08023483365949119
$ws.Range("F18").Value = 0.0131578947368421
$ws.Range("H18").Value = 0.1710526315789474
$ws.Range("I18").Value = 0.07236842105263158
$ws.Range("J18").Value = 0.4210526315789473
$ws.Range("K18").Value = 0.1184210526315789
$ws.Range("M18").Value = 0.006578947368421052
$ws.Range("O18").Value = 0.09210526315789473
$ws.Range("S18").Value = 0.1052631578947368
$ws.Range("F19").Value = 0.01812450748620961
$ws.Range("H19").Value = 0.2040977147360126
$ws.Range("I19").Value = 0.09771473601260836
$ws.Range("J19").Value = 0.3798266351457841
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.02048857368006304
$ws.Range("N19").Value = 0.0007880220646178094
$ws.Range("O19").Value = 0.0677698975571316
$ws.Range("S19").Value = 0.1000788022064618
